$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-17 Friday" "2025-10-18 Saturday"

Replace-Text "23÷8=" "57÷4="
Replace-Text "50÷9=" "55÷8="
Replace-Text "74÷3=" "96÷7="
Replace-Text "20÷3=" "46÷8="
Replace-Text "29÷9=" "18÷6="
Replace-Text "40÷2=" "30÷4="
Replace-Text "33÷9=" "89÷5="
Replace-Text "34÷9=" "49÷9="
Replace-Text "35÷6=" "45÷7="
Replace-Text "32÷9=" "60÷5="
Replace-Text "16÷2=" "21÷6="
Replace-Text "37÷3=" "21÷2="
Replace-Text "46÷4=" "62÷6="
Replace-Text "49÷8=" "91÷6="
Replace-Text "50÷8=" "32÷2="
Replace-Text "46÷9=" "15÷7="
Replace-Text "23÷9=" "99÷9="
Replace-Text "13÷2=" "70÷2="
Replace-Text "44÷6=" "90÷4="
Replace-Text "51÷7=" "52÷8="
Replace-Text "98÷2=" "78÷8="
Replace-Text "75÷2=" "73÷3="
Replace-Text "26÷4=" "12÷5="
Replace-Text "78÷4=" "64÷5="
Replace-Text "80÷4=" "59÷4="
